# Global Software Control.docx — update the explanation of how the
# application server dispatches incoming requests: the server now
# hands the redirect off to a dedicated thread (to serve every
# connected user concurrently) before forwarding to the dispatching
# endpoints.

$d = $word.ActiveDocument

$old = " server, che le re-indirizza sugli appositi endpoint di "
$new = " server, che tramite un thread dedicato, in modo da garantire " + `
       "un" + [char]0x2019 + "interazione concorrente con tutti gli utenti connessi, " + `
       "li re-indirizza sugli appositi endpoint di "

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
